$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": status text changed from "Ready for handoff" to
# "Handed back: in sync with en-US" for both language columns (E,F) on both
# data rows (2,3). Columns E & F also get wider to fit the longer text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.083333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 29.083333333333336

# ---------------------------------------------------------------------------
# Helper values shared by both per-language report sheets.
# ---------------------------------------------------------------------------
$mdUrl668 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md"
$mdUrl99b = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md"
$md668 = "668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md"
$md99b = "99b90781-9224-4582-ba7b-4fe81cf19a3c.md"

$hyperlinkColor = 15570276   # cornflower blue (FF6495ED) used by the workbook's "HyperLink" style

# ---------------------------------------------------------------------------
# Sheet "zh-cn": the handback run resolved the "Latest Target File" (I) and
# "Latest Handback File" (J) columns, plus the "Latest Handback DateTime" (K).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = $md668
$wsZh.Range("J2").Value = "668e1d8b-ca32-4470-bc2d-b3a2537e67a4.3b9d6a9edaff5c3fdef856fe09cef804ea4d6893.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 11:07:54"

$wsZh.Range("I3").Value = $md99b
$wsZh.Range("J3").Value = "99b90781-9224-4582-ba7b-4fe81cf19a3c.98c20a0cba445a2e320b9a43ed030525c5f87c69.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 11:07:54"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl668, [Type]::Missing, [Type]::Missing, $md668)
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = $hyperlinkColor

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl99b, [Type]::Missing, [Type]::Missing, $md99b)
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = $hyperlinkColor

$wsZh.Columns.Item(3).ColumnWidth = 29.083333333333336
$wsZh.Columns.Item(9).ColumnWidth = 39.08333333333333
$wsZh.Columns.Item(10).ColumnWidth = 39.08333333333333

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape of update, but the handback run for this
# language landed at a different timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = $md668
$wsDe.Range("J2").Value = "668e1d8b-ca32-4470-bc2d-b3a2537e67a4.3b9d6a9edaff5c3fdef856fe09cef804ea4d6893.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 11:08:10"

$wsDe.Range("I3").Value = $md99b
$wsDe.Range("J3").Value = "99b90781-9224-4582-ba7b-4fe81cf19a3c.98c20a0cba445a2e320b9a43ed030525c5f87c69.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 11:08:10"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl668, [Type]::Missing, [Type]::Missing, $md668)
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = $hyperlinkColor

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl99b, [Type]::Missing, [Type]::Missing, $md99b)
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = $hyperlinkColor

$wsDe.Columns.Item(3).ColumnWidth = 29.083333333333336
$wsDe.Columns.Item(9).ColumnWidth = 39.08333333333333
$wsDe.Columns.Item(10).ColumnWidth = 39.08333333333333
